# Fill in the "year" column (B) for the existing rows (2-93) with the
# sequential years 1929-2020. This also replaces the old text value
# "2019 or 2020" that lived in B93 with the numeric year 2020.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$year = 1929
for ($row = 2; $row -le 93; $row++) {
    $ws.Cells.Item($row, 2).Value = $year
    $year = $year + 1
}

# Append three new data rows (94-96) at the bottom of the table.
$ws.Cells.Item(94, 1).Value = 93
$ws.Cells.Item(94, 2).Value = 1945
$ws.Cells.Item(94, 3).Value = 18
$ws.Cells.Item(94, 4).Value = 2

$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = 1955
$ws.Cells.Item(95, 3).Value = 28
$ws.Cells.Item(95, 4).Value = 2

$ws.Cells.Item(96, 1).Value = 95
$ws.Cells.Item(96, 2).Value = 2019
$ws.Cells.Item(96, 3).Value = 92
$ws.Cells.Item(96, 4).Value = 2

# Scroll / select to match where the author left the sheet after editing.
[void]$ws.Range("C97").Select()
